$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.066.21"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.643.63"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "214.68"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "0.5101"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.2566"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "0.06359"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").Value = "4.289"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "1.648.40"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "0.5445"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "64.37"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0₅7718"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "26.046.30"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "198.39"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "4.427"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "9.929"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").Value = "6.042"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "1.865"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "141.01"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").Value = "0.1196"
$ws.Range("E26").Value = "  +5.32%  "
$ws.Range("D27").Value = "6.819"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "15.63"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").Value = "0.04861"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "3.256"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "3.172"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "1.527"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").Value = "2.360"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "0.8994"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "2.581"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "1.142.34"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "0.5470"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "0.01563"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "1.000"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "2.532"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").Value = "0.0₈129"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("D43").Value = "0.8110"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "99.36"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "5.391"
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("D46").Value = "1.779.06"
$ws.Range("D47").Value = "0.4526"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "54.92"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "0.9987"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "0.05057"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.39%  "
